$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = 0
$ws.Range("H12").Value = 1749.875
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 1000
$ws.Range("M12").Value = -830
$ws.Range("H33").Value = 180.4
$ws.Range("I33").Value = 163.25
$ws.Range("K33").Value = 163.25
$ws.Range("M33").Value = 65.75
$ws.Range("H39").Value = 1614.25
$ws.Range("I39").Value = 2073.4
$ws.Range("J39").Value = 849
$ws.Range("K39").Value = 6220.200000000001
$ws.Range("L39").Value = 2547
$ws.Range("M39").Value = -5924.200000000001
$ws.Range("N39").Value = -3139
$ws.Range("H64").Value = 3899.9285
$ws.Range("I64").Value = 3899.9285
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3899.9285
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -3651.9285
$ws.Range("H67").Value = 3899.9285
$ws.Range("I67").Value = 3899.9285
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3899.9285
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -3041.9285
$ws.Range("H74").Value = 4483.4
$ws.Range("I74").Value = 3417
$ws.Range("K74").Value = 3417
$ws.Range("M74").Value = -2481
$ws.Range("H77").Value = 4483.4
$ws.Range("I77").Value = 3417
$ws.Range("K77").Value = 17085
$ws.Range("M77").Value = -12405
$ws.Range("H92").Value = 361.65
$ws.Range("I92").Value = 361.65
$ws.Range("K92").Value = 361.65
$ws.Range("M92").Value = 886.35
$ws.Range("H113").Value = 45093.043
$ws.Range("I113").Value = 102332.6
$ws.Range("K113").Value = 102332.6
$ws.Range("M113").Value = -99078.60000000001
$ws.Range("H137").Value = 1093.6923
$ws.Range("I137").Value = 1023.5455
$ws.Range("K137").Value = 3070.6365
$ws.Range("M137").Value = -520.6364999999996

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6151.7534
$ws.Range("I32").Value = 2486.164
$ws.Range("K32").Value = 2486.164
$ws.Range("M32").Value = -2199.164
$ws.Range("H122").Value = 2568.7778
$ws.Range("I122").Value = 2564.5
$ws.Range("K122").Value = 7693.5
$ws.Range("M122").Value = -5243.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1752.6086
$ws.Range("I86").Value = 1415.0714
$ws.Range("K86").Value = 1415.0714
$ws.Range("M86").Value = -292.0714
$ws.Range("H89").Value = 1752.6086
$ws.Range("I89").Value = 1415.0714
$ws.Range("K89").Value = 7075.357
$ws.Range("M89").Value = -1459.357
$ws.Range("H134").Value = 2976.6155
$ws.Range("I134").Value = 924.75
$ws.Range("J134").Value = 3888.5557
$ws.Range("K134").Value = 2774.25
$ws.Range("L134").Value = 11665.6671
$ws.Range("M134").Value = -239.25
$ws.Range("N134").Value = -16735.6671

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3388.7144
$ws.Range("I16").Value = 5669.6665
$ws.Range("J16").Value = 1678
$ws.Range("K16").Value = 5669.6665
$ws.Range("L16").Value = 1678
$ws.Range("M16").Value = -5382.6665
$ws.Range("N16").Value = -2252
$ws.Range("H31").Value = 37739.08
$ws.Range("I31").Value = 33951.375
$ws.Range("K31").Value = 33951.375
$ws.Range("M31").Value = -33656.375
$ws.Range("H34").Value = 37739.08
$ws.Range("I34").Value = 33951.375
$ws.Range("K34").Value = 33951.375
$ws.Range("M34").Value = -33749.375
$ws.Range("H105").Value = 1640.3334
$ws.Range("I105").Value = 910
$ws.Range("K105").Value = 910
$ws.Range("M105").Value = 837
$ws.Range("H113").Value = 3388.7144
$ws.Range("I113").Value = 5669.6665
$ws.Range("J113").Value = 1678
$ws.Range("K113").Value = 5669.6665
$ws.Range("L113").Value = 1678
$ws.Range("M113").Value = -3499.6665
$ws.Range("N113").Value = -6018
$ws.Range("H132").Value = 3852.158
$ws.Range("I132").Value = 3791.3333
$ws.Range("K132").Value = 11373.9999
$ws.Range("M132").Value = -8843.999899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 333377300
$ws.Range("J37").Value = 333377300
$ws.Range("L37").Value = 1000131900
$ws.Range("N37").Value = -1000132124
$ws.Range("H116").Value = 6250
$ws.Range("I116").Value = 2500
$ws.Range("J116").Value = 10000
$ws.Range("K116").Value = 7500
$ws.Range("L116").Value = 30000
$ws.Range("M116").Value = -4058
$ws.Range("N116").Value = -36884

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8147.625
$ws.Range("I70").Value = 7951.3076
$ws.Range("J70").Value = 8998.333000000001
$ws.Range("K70").Value = 7951.3076
$ws.Range("L70").Value = 8998.333000000001
$ws.Range("M70").Value = -7681.3076
$ws.Range("N70").Value = -9538.333000000001
$ws.Range("H73").Value = 8147.625
$ws.Range("I73").Value = 7951.3076
$ws.Range("J73").Value = 8998.333000000001
$ws.Range("K73").Value = 7951.3076
$ws.Range("L73").Value = 8998.333000000001
$ws.Range("M73").Value = -7015.3076
$ws.Range("N73").Value = -10870.333
$ws.Range("H129").Value = 40390
$ws.Range("J129").Value = 30780
$ws.Range("L129").Value = 30780
$ws.Range("N129").Value = -40780
$ws.Range("H132").Value = 8695.571
$ws.Range("I132").Value = 7869.857
$ws.Range("K132").Value = 23609.571
$ws.Range("M132").Value = -21079.571

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1068.1428
$ws.Range("I22").Value = 1036.75
$ws.Range("K22").Value = 1036.75
$ws.Range("M22").Value = -741.75
$ws.Range("H27").Value = 1068.1428
$ws.Range("I27").Value = 1036.75
$ws.Range("K27").Value = 1036.75
$ws.Range("M27").Value = -929.75
$ws.Range("H55").Value = 327.73077
$ws.Range("I55").Value = 353.6111
$ws.Range("J55").Value = 269.5
$ws.Range("K55").Value = 353.6111
$ws.Range("L55").Value = 269.5
$ws.Range("M55").Value = -180.6111
$ws.Range("N55").Value = -615.5
$ws.Range("H61").Value = 90743.2
$ws.Range("I61").Value = 92847.37
$ws.Range("J61").Value = 84956.75
$ws.Range("K61").Value = 92847.37
$ws.Range("L61").Value = 84956.75
$ws.Range("M61").Value = -92645.37
$ws.Range("N61").Value = -85360.75
$ws.Range("H113").Value = 90743.2
$ws.Range("I113").Value = 92847.37
$ws.Range("J113").Value = 84956.75
$ws.Range("K113").Value = 92847.37
$ws.Range("L113").Value = 84956.75
$ws.Range("M113").Value = -90677.37
$ws.Range("N113").Value = -89296.75
$ws.Range("H132").Value = 3477.111
$ws.Range("I132").Value = 3386.9644
$ws.Range("J132").Value = 3792.625
$ws.Range("K132").Value = 10160.8932
$ws.Range("L132").Value = 11377.875
$ws.Range("M132").Value = -7630.893199999999
$ws.Range("N132").Value = -16437.875
$ws.Range("H133").Value = 113331.664
$ws.Range("J133").Value = 113331.664
$ws.Range("L133").Value = 113331.664
$ws.Range("N133").Value = -118391.664

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 12
$ws.Range("I26").Value = 12
$ws.Range("K26").Value = 12
$ws.Range("M26").Value = 281
$ws.Range("H81").Value = 1578.3125
$ws.Range("I81").Value = 1703.7858
$ws.Range("K81").Value = 3407.5716
$ws.Range("M81").Value = -2346.5716
$ws.Range("H84").Value = 1578.3125
$ws.Range("I84").Value = 1703.7858
$ws.Range("K84").Value = 17037.858
$ws.Range("M84").Value = -11733.858
$ws.Range("H113").Value = 411
$ws.Range("I113").Value = 122
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 366
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = 1804
$ws.Range("N113").Value = -6440
$ws.Range("H133").Value = 69108.39999999999
$ws.Range("J133").Value = 69108.39999999999
$ws.Range("L133").Value = 69108.39999999999
$ws.Range("N133").Value = -79228.39999999999

Write-Host "Applied market-price refresh edits"